# Weekly Timesheet (sheet1) updates
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rename client/employee names (B column, rows 2-6)
$ws1.Range("B2").Value = "Hunter"
$ws1.Range("B3").Value = "Tubergen"
$ws1.Range("B4").Value = "Field"
$ws1.Range("B5").Value = "Bottomley"
$ws1.Range("B6").Value = "Zygmunt"

# Update hours/rate/total for rows 2-5 (8 -> 9 hrs, rate 0 -> 92, total 0 -> 828)
$ws1.Range("C2").Value = 9
$ws1.Range("E2").Value = 92
$ws1.Range("F2").Value = 828

$ws1.Range("C3").Value = 9
$ws1.Range("E3").Value = 92
$ws1.Range("F3").Value = 828

$ws1.Range("C4").Value = 9
$ws1.Range("E4").Value = 92
$ws1.Range("F4").Value = 828

$ws1.Range("C5").Value = 9
$ws1.Range("E5").Value = 92
$ws1.Range("F5").Value = 828

# Row 6: hours 8 -> 4, rate -> 92, total -> 368 (remainder moved to new OT row)
$ws1.Range("C6").Value = 4
$ws1.Range("E6").Value = 92
$ws1.Range("F6").Value = 368

# Insert a new row 7 (duplicate of row 6's formatting) for the OT entry
$ws1.Rows("6").Copy()
$ws1.Rows("7").Insert()
$ws1.Range("C7").Value = 5
$ws1.Range("D7").Value = "OT"
$ws1.Range("F7").Value = 690

# SUBTOTAL row (was row 8, now row 9 after the insert): hours 40 -> 45, text, total -> 4370
$ws1.Range("C9").Value = 45
$ws1.Range("D9").Value = "Reg: 40 / OT: 5"
$ws1.Range("F9").Value = 4370

# HOURLY SUBTOTAL row (was row 11, now row 12): total 0 -> 4370
$ws1.Range("F12").Value = 4370

# GRAND TOTAL row (was row 13, now row 14): total 0 -> 4370
$ws1.Range("F14").Value = 4370

# Jason Schema (sheet2) updates
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B2").Value = "emp_yde33znx"
$ws2.Range("B3").Value = "emp_yde33znx"
$ws2.Range("B4").Value = "emp_yde33znx"
$ws2.Range("B5").Value = "emp_yde33znx"
$ws2.Range("B6").Value = "emp_yde33znx"

$ws2.Range("D2").Value = "Hunter"
$ws2.Range("D3").Value = "Tubergen"
$ws2.Range("D4").Value = "Field"
$ws2.Range("D5").Value = "Bottomley"
$ws2.Range("D6").Value = "Zygmunt"

$ws2.Range("E2").Value = 9
$ws2.Range("F2").Value = 92
$ws2.Range("G2").Value = 828

$ws2.Range("E3").Value = 9
$ws2.Range("F3").Value = 92
$ws2.Range("G3").Value = 828

$ws2.Range("E4").Value = 9
$ws2.Range("F4").Value = 92
$ws2.Range("G4").Value = 828

$ws2.Range("E5").Value = 9
$ws2.Range("F5").Value = 92
$ws2.Range("G5").Value = 828

$ws2.Range("E6").Value = 4
$ws2.Range("F6").Value = 92
$ws2.Range("G6").Value = 368

# New row 7 for the OT entry (duplicate row 6 formatting, then adjust)
$ws2.Rows("6").Copy()
$ws2.Rows("7").Insert()
$ws2.Range("E7").Value = 5
$ws2.Range("G7").Value = 690
$ws2.Range("H7").Value = "OT"
